$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $r = $cell.Range
    $r.MoveEnd(1, -1) | Out-Null
    $r.Text = $newText
}

# Simple value replacements (by row, 1-indexed)
Set-CellText $t 1 1 "0M"
Set-CellText $t 2 1 "0M"
Set-CellText $t 3 1 "0M"
Set-CellText $t 4 1 "7985"
Set-CellText $t 7 1 "0.02231"
Set-CellText $t 8 1 "0.01082"
Set-CellText $t 12 1 "73.10511"

# Collapse the multi-run tab-separated cells down to the single prior summary value
Set-CellText $t 44 1 "90.31"
Set-CellText $t 45 1 "73.11"
Set-CellText $t 46 1 "754"
